$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.900.17"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "'3.032.00"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'586.14"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'149.37"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("D9").Value = "'3.033.05"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").Value = "'35.33"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "'0.122"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").Value = "'3.535.66"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'62.845.96"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'3.032.53"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'468.75"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "'14.05"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'0.692"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").Value = "'80.91"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'10.41"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.26"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "'2.16"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'27.74"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  -4.17%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'0.0₃0806"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "'5.79"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "'2.15"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Value = "'50.35"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "'9.02"
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").Value = "'2.97"
$ws.Range("E41").Value = "  -10.05%  "
$ws.Range("D42").Value = "'426.88"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "'0.281"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'2.805.09"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'0.0356"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'37.79"
$ws.Range("E47").Value = "  -9.07%  "
$ws.Range("D48").Value = "'129.43"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D50").Value = "'24.43"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("E51").Value = "  -0.36%  "
